$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated category labels and values (a new "Abnormal test result" row was
# inserted at row 2, shifting the remaining categories down by one row, and
# all numeric values were refreshed per the corrected POC pCO2 extraction).

$ws.Range("A2").Value = "Abnormal test result"
$ws.Range("B2").Value = 1.6
$ws.Range("C2").Value = 1.5
$ws.Range("D2").Value = 2.6

$ws.Range("A3").Value = "Diseases (patient-stated)"
$ws.Range("B3").Value = 4.9
$ws.Range("C3").Value = 6.2
$ws.Range("D3").Value = 6.6

$ws.Range("A4").Value = "Injuries & adverse effects"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 12.1
$ws.Range("D4").Value = 5.8

$ws.Range("A5").Value = "Other"
$ws.Range("B5").Value = 6.5
$ws.Range("C5").Value = 9.199999999999999
$ws.Range("D5").Value = 6.4

$ws.Range("A6").Value = "Symptom – Circulatory"
$ws.Range("B6").Value = 7.7
$ws.Range("C6").Value = 8.6
$ws.Range("D6").Value = 9.9

$ws.Range("A7").Value = "Symptom – Digestive"
$ws.Range("B7").Value = 9.699999999999999
$ws.Range("C7").Value = 12.8
$ws.Range("D7").Value = 15.6

$ws.Range("A8").Value = "Symptom – General"
$ws.Range("B8").Value = 4.7
$ws.Range("C8").Value = 5.3
$ws.Range("D8").Value = 7

$ws.Range("A9").Value = "Symptom – Genitourinary"
$ws.Range("B9").Value = 4.8
$ws.Range("C9").Value = 6.6
$ws.Range("D9").Value = 5.7

$ws.Range("A10").Value = "Symptom – Nervous"
$ws.Range("B10").Value = 11.3
$ws.Range("C10").Value = 9.800000000000001
$ws.Range("D10").Value = 14.5

$ws.Range("A11").Value = "Symptom – Respiratory"
$ws.Range("B11").Value = 38.4
$ws.Range("C11").Value = 24.9
$ws.Range("D11").Value = 23.8

$ws.Range("A12").Value = "Symptom – Skin/Hair/Nails"
$ws.Range("B12").Value = 2.4
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 2.1
